$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column width updates (C, D, F, G, H change; A, B, E stay the same).
# Excel's COM ColumnWidth is expressed in "characters" and is offset from the
# stored OOXML column width by the standard ~0.8333333 character padding, so
# we back that out here to land on the exact target widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 39.1666667   # -> 40
$ws.Columns.Item(4).ColumnWidth = 43.1666667   # -> 44
$ws.Columns.Item(6).ColumnWidth = 16.1666667   # -> 17
$ws.Columns.Item(7).ColumnWidth = 15.1666667   # -> 16
$ws.Columns.Item(8).ColumnWidth = 44.1666667   # -> 45

# ---------------------------------------------------------------------------
# Row 2 — new scraped opportunity record.
# Opportunity IDs are numeric-looking but must stay text cells (as in the
# rest of the sheet). Entering them as a quoted formula and then pasting
# back as values keeps the cell's text type without leaving the "General"
# style touched (a plain numeric literal, or an apostrophe text-prefix,
# would otherwise coerce the cell to a Number or tag it with quotePrefix).
# ---------------------------------------------------------------------------
$ws.Range("A2").Formula = "=""1330604"""
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1330604"
$ws.Range("C2").Value = "[EXP] People Data Specialist Intern"
$ws.Range("D2").Value = "Fritz-Erler-Straße 5, 53113 Bonn, Germany"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "2 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "DHL Group"

# Highlight the PREMIUM cell for this row with a solid yellow fill.
$ws.Range("E2").Interior.Color = 65535

# ---------------------------------------------------------------------------
# Row 3 — new scraped opportunity record.
# ---------------------------------------------------------------------------
$ws.Range("A3").Formula = "=""1328614"""
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)

$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328614"
$ws.Range("C3").Value = "Field Service Engineer [EU Preferred]"
$ws.Range("D3").Value = "Madrid, Spain"
$ws.Range("F3").Value = "126 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Mitsubishi Power Europe Sucursal en España"
